$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Inscritos (E) 10 -> 12
$ws.Range("E7").Value = 12

# Row 13: Inscritos (E) 1 -> 2
$ws.Range("E13").Value = 2

# Row 17: Pagos (F) 71 -> 72, Inscrições homologadas (H) 103 -> 104
$ws.Range("F17").Value = 72
$ws.Range("H17").Value = 104

# Row 36: Inscritos (E) 116 -> 118
$ws.Range("E36").Value = 118

# Row 43: Inscritos (E) 30 -> 31
$ws.Range("E43").Value = 31

# Row 58: Inscritos (E) 4 -> 5, Pagos (F) 3 -> 4, Inscrições homologadas (H) 3 -> 4
$ws.Range("E58").Value = 5
$ws.Range("F58").Value = 4
$ws.Range("H58").Value = 4

# Row 62: Inscritos (E) 51 -> 52
$ws.Range("E62").Value = 52

# Row 67: Inscritos (E) 43 -> 44
$ws.Range("E67").Value = 44

# Row 77: Inscritos (E) 62 -> 63
$ws.Range("E77").Value = 63
